$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "RFID Number" header in D1 to "RFIDNumber" (remove the space)
$ws.Range("D1").Value = "RFIDNumber"

# Selection now only covers D1 (previously D1:D11)
$ws.Range("D1").Select()
